# Backup/Payments.xlsx update:
#  - G12 (IsDeleted) flips from 0 -> 1
#  - New rows 13-18 appended with payment records
#  - dimension grows to A1:G18 (handled automatically by the engine)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing row 12: mark as deleted
$ws.Range("G12").Value = 1

# New rows appended after row 12
$data = @(
    @(12, 19, 6, 10, 0,     "2025-03-27 17:47:33", 0),
    @(13, 18, 7, 11, 0,     "2025-03-27 17:50:33", 0),
    @(14, 20, 8, 11, 0,     "2025-03-28 18:18:25", 0),
    @(15, 20, 8, 11, 10000, "2025-03-28 18:19:38", 0),
    @(16, 20, 9, 11, 0,     "2025-03-28 18:41:45", 0),
    @(17, 20, 9, 11, 9000,  "2025-03-28 18:42:51", 0)
)

$rowIndex = 13
foreach ($record in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $record[0]
    $ws.Cells.Item($rowIndex, 2).Value = $record[1]
    $ws.Cells.Item($rowIndex, 3).Value = $record[2]
    $ws.Cells.Item($rowIndex, 4).Value = $record[3]
    $ws.Cells.Item($rowIndex, 5).Value = $record[4]
    $ws.Cells.Item($rowIndex, 6).Value = $record[5]
    $ws.Cells.Item($rowIndex, 7).Value = $record[6]
    $rowIndex++
}
